{"js": "// Apply the \"UserData template\" edit: prefix every placeholder paragraph's\n// text with \"$\" and switch the placeholder syntax from \"{name}\" style to\n// \"{.name}\" (dot-prefixed) style; the closing \"{/records}\" tag becomes\n// \"{records#}\" (still \"$\"-prefixed); two trailing helper paragraphs\n// (\"{#name_tag}\" and \"{#html_field}\") are removed entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Map of the exact current paragraph text -> the new paragraph text.\n// insertText(..., \"Replace\") rewrites the paragraph's text while keeping\n// the formatting of the run(s) it replaces.\nconst replacements = {\n  \"{#records}\": \"${#records}\",\n  \"{name}\": \"${.name}\",\n  \"{user_name}\": \"${.user_name}\",\n  \"{email}\": \"${.email}\",\n  \"{gender}\": \"${.gender}\",\n  \"{u_national_holiday_country}\": \"${.u_national_holiday_country}\",\n  \"{/records}\": \"${records#}\",\n};\n\n// Paragraphs whose text content should cause the whole paragraph to be\n// deleted (the trailing helper placeholders are no longer needed).\nconst toDelete = [\"{#name_tag}\", \"{#html_field}\"];\n\nfor (const para of items) {\n  const text = para.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, text)) {\n    para.insertText(replacements[text], \"Replace\");\n  } else if (toDelete.includes(text)) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the \"UserData template\" edit: prefix every placeholder paragraph's\n# text with \"$\" and switch the placeholder syntax from \"{name}\" style to\n# \"{.name}\" (dot-prefixed) style; the closing \"{/records}\" tag becomes\n# \"{records#}\" (still \"$\"-prefixed); two trailing helper paragraphs\n# (\"{#name_tag}\" and \"{#html_field}\") are removed entirely.\n\n$d = $word.ActiveDocument\n\n# Delete the two trailing helper paragraphs first (from the end, so the\n# earlier paragraph indices/ranges used below stay valid).\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n    if ($t -eq \"{#name_tag}\" -or $t -eq \"{#html_field}\") {\n        $p.Range.Delete()\n    }\n}\n\n# Text replacements, using literal Find/Replace (no wildcards) so the\n# curly braces in the search/replacement strings are treated as plain text.\n$pairs = @(\n    @(\"{#records}\", \"`${#records}\"),\n    @(\"{name}\", \"`${.name}\"),\n    @(\"{user_name}\", \"`${.user_name}\"),\n    @(\"{email}\", \"`${.email}\"),\n    @(\"{gender}\", \"`${.gender}\"),\n    @(\"{u_national_holiday_country}\", \"`${.u_national_holiday_country}\"),\n    @(\"{/records}\", \"`${records#}\")\n)\n\nforeach ($pair in $pairs) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
